# Updates betexplorer scrape sheet:
#  - three pairs of adjacent rows had their match data (columns F:V) swapped
#    (row order bug fixed upstream), identity columns A:E untouched
#  - three brand-new match rows appended at the end (96, 97, 98)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchData($rowA, $rowB) {
    for ($c = 6; $c -le 22; $c++) {
        $valA = $ws.Cells.Item($rowA, $c).Value2
        $valB = $ws.Cells.Item($rowB, $c).Value2
        $ws.Cells.Item($rowA, $c).Value2 = $valB
        $ws.Cells.Item($rowB, $c).Value2 = $valA
    }
}

# --- swap mis-ordered row pairs -------------------------------------------
Swap-MatchData 15 16
Swap-MatchData 50 51
Swap-MatchData 53 54

# --- append new rows 96-98, cloning the formatting of the last data row ---
$ws.Range("A95:V95").Copy()
$ws.Range("A96:V98").PasteSpecial(-4122)  # xlPasteFormats

function Set-MatchRow {
    param(
        $row, $indice, $dataPartida,
        $home, $homeGols, $away, $awayGols,
        $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
        $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
        $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
        $url
    )
    $ws.Cells.Item($row, 1).Value2 = $indice
    $ws.Cells.Item($row, 2).Value2 = "turkey"
    $ws.Cells.Item($row, 3).Value2 = "1-lig"
    $ws.Cells.Item($row, 4).Value2 = "2023-2024"
    $ws.Cells.Item($row, 5).Value2 = $dataPartida
    $ws.Cells.Item($row, 6).Value2 = $home
    $ws.Cells.Item($row, 7).Value2 = $homeGols
    $ws.Cells.Item($row, 8).Value2 = $away
    $ws.Cells.Item($row, 9).Value2 = $awayGols
    $ws.Cells.Item($row, 10).Value2 = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value2 = $homeOpenDt
    $ws.Cells.Item($row, 12).Value2 = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value2 = $homeCloseDt
    $ws.Cells.Item($row, 14).Value2 = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value2 = $drawOpenDt
    $ws.Cells.Item($row, 16).Value2 = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value2 = $drawCloseDt
    $ws.Cells.Item($row, 18).Value2 = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value2 = $awayOpenDt
    $ws.Cells.Item($row, 20).Value2 = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value2 = $awayCloseDt
    $ws.Cells.Item($row, 22).Value2 = $url
}

Set-MatchRow 96 95 45235.47916666666 `
    "Bodrumspor" 2 "Manisa FK" 0 `
    1.95 "29/10/2023 11:42" 1.97 "04/11/2023 20:42" `
    3.43 "29/10/2023 11:42" 3.44 "04/11/2023 20:42" `
    3.98 "29/10/2023 11:42" 4.03 "04/11/2023 20:42" `
    "https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-manisa-fk/WIu9cXQo/"

Set-MatchRow 97 96 45235.47916666666 `
    "Bandirmaspor" 2 "Altay" 0 `
    1.53 "29/10/2023 17:13" 1.33 "05/11/2023 11:28" `
    4.25 "29/10/2023 17:13" 5.29 "05/11/2023 11:28" `
    5.78 "29/10/2023 17:13" 9.06 "05/11/2023 11:28" `
    "https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-altay/4WwHeBeb/"

Set-MatchRow 98 97 45235.58333333334 `
    "Giresunspor" 1 "Erzurumspor" 0 `
    2.79 "29/10/2023 14:12" 3.35 "05/11/2023 13:54" `
    3.08 "29/10/2023 14:12" 3.32 "05/11/2023 13:54" `
    2.69 "29/10/2023 14:12" 2.25 "05/11/2023 13:54" `
    "https://www.betexplorer.com/football/turkey/1-lig/giresunspor-erzurumspor-fk/vkWTh9uH/"
